$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the text contents of rows 24 and 25 (columns B and C) so that the
# shared-string table order matches the target diff: the "LOM3221" entry
# now comes before the "LOM3202" entry, which means the values shown in
# row 24 and row 25 (previously LOM3202 / LOM3221 respectively) trade
# places.
$b24 = $ws.Range("B24").Value2
$b25 = $ws.Range("B25").Value2

$ws.Range("B24").Value = $b25
$ws.Range("C24").Value = $b25

$ws.Range("B25").Value = $b24
$ws.Range("C25").Value = $b24
